$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New "age group" header column inserted in front of the existing table
#    (row 1, column A). The rest of row 1's headers (B1:G1) stay put.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "age group"

# ---------------------------------------------------------------------------
# 2. Column width tweaks: column A becomes wider (to fit the new "age group"
#    labels) and column F becomes wider too (for the new header text used
#    further down in the sheet). Excel stores column widths in pixel-
#    quantised character units, so we dial in the ColumnWidth that produces
#    the closest on-disk width to the target.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 19.166666666666668
$ws.Columns.Item(6).ColumnWidth = 15.166666666666666

# ---------------------------------------------------------------------------
# 3. New block below the original table (rows 25-38): a second, similar
#    "age group" table that derives a fitted/pop-weighted/desired
#    hospitalisation propensity ("h") per age group.
# ---------------------------------------------------------------------------

# --- header row (25) --------------------------------------------------
$ws.Range("A25").Value = "age group"
$ws.Range("B25").Value = "n0 individuals"
$ws.Range("C25").Value = "fraction"
$ws.Range("D25").Value = "fitted h"
$ws.Range("E25").Value = "multiplier"
$ws.Range("F25").Value = "pop weighted h"
$ws.Range("G25").Value = "desired h"
$ws.Range("H25").Value = "rounded h (%)"

# --- per age-group data rows (26-34) -----------------------------------
$ageGroups = @("0-10", "10-20", "20-30", "30-40", "40-50", "50-60", "60-70", "70-80", "80-inf")
$nIndividuals = @(1305219, 1298970, 1395385, 1498535, 1524152, 1601891, 1347696, 908725, 658753)
$fittedH = @(0.015, 0.02, 0.03, 0.03, 0.03, 0.06, 0.14, 0.3, 0.76)

for ($i = 0; $i -lt 9; $i++) {
    $r = 26 + $i
    $ws.Range("A$r").Value = $ageGroups[$i]
    $ws.Range("B$r").Value = $nIndividuals[$i]
    $ws.Range("C$r").Formula = "=B$r/`$B`$35"
    $ws.Range("D$r").Value = $fittedH[$i]
    $ws.Range("F$r").Formula = "=D$r*`$E`$26*C$r"
    $ws.Range("G$r").Formula = "=D$r*`$E`$26"
    $ws.Range("H$r").Formula = "=ROUND(G$r,3)*100"
}

# E26 holds the (literal) population-weighted multiplier used by every row
# via the absolute reference $E$26.
$ws.Range("E26").Value = 0.746784953353961

# --- totals row (35) -----------------------------------------------------
$ws.Range("A35").Value = "total"
$ws.Range("B35").Formula = "=SUM(B26:B34)"
$ws.Range("C35").Formula = "=B35/`$B`$35"
$ws.Range("F35").Formula = "=SUM(F26:F34)"
$ws.Range("F35").Font.Bold = $true

# --- desired pop weighted h (row 38) --------------------------------------
$ws.Range("A38").Value = "desired pop weighted h"
$ws.Range("B38").Value = 0.08
$ws.Range("B38").Font.Bold = $true

# ---------------------------------------------------------------------------
# 4. Restore the active selection to match the author's final cursor
#    position.
# ---------------------------------------------------------------------------
$ws.Range("F19").Select()
